$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "nuts"
$ws.Range("A5").Value = "honey"
$ws.Range("A6").Value = "fruits"
$ws.Range("A4").Value = "breads"

$ws.Range("A4").Select() | Out-Null
